# Update "想去人数" (F column) values across sheets to reflect the latest
# generated output (commit: Update gh-pages to output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 150
$ws.Range("F4").Value = 812
$ws.Range("F6").Value = 690
$ws.Range("F7").Value = 1251
$ws.Range("F9").Value = 854
$ws.Range("F10").Value = 711
$ws.Range("F11").Value = 267
$ws.Range("F13").Value = 380
$ws.Range("F15").Value = 1005
$ws.Range("F16").Value = 11223
$ws.Range("F17").Value = 645
$ws.Range("F18").Value = 52
$ws.Range("F22").Value = 284
$ws.Range("F23").Value = 1791
$ws.Range("F24").Value = 30
$ws.Range("F25").Value = 293
$ws.Range("F26").Value = 494
$ws.Range("F27").Value = 190
$ws.Range("F29").Value = 288
$ws.Range("F30").Value = 200
$ws.Range("F32").Value = 78
$ws.Range("F35").Value = 183
$ws.Range("F36").Value = 202
$ws.Range("F37").Value = 301

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 147
$ws.Range("F12").Value = 87
$ws.Range("F16").Value = 321

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 830

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 830
$ws.Range("F5").Value = 150
$ws.Range("F6").Value = 812
$ws.Range("F9").Value = 690
$ws.Range("F10").Value = 1251
$ws.Range("F13").Value = 147
$ws.Range("F14").Value = 854
$ws.Range("F15").Value = 711
$ws.Range("F16").Value = 267
$ws.Range("F18").Value = 1005
$ws.Range("F19").Value = 11223
$ws.Range("F21").Value = 645
$ws.Range("F22").Value = 52
$ws.Range("F24").Value = 284
$ws.Range("F25").Value = 1791
$ws.Range("F26").Value = 494
$ws.Range("F27").Value = 190
$ws.Range("F28").Value = 87
$ws.Range("F29").Value = 87
$ws.Range("F33").Value = 321
$ws.Range("F34").Value = 288
$ws.Range("F36").Value = 200
$ws.Range("F38").Value = 78
$ws.Range("F42").Value = 183
$ws.Range("F45").Value = 202
$ws.Range("F46").Value = 305
